# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the newly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$wsExhibition.Range("F3").Value  = 2223
$wsExhibition.Range("F5").Value  = 13339
$wsExhibition.Range("F8").Value  = 521
$wsExhibition.Range("F10").Value = 1195
$wsExhibition.Range("F11").Value = 1003
$wsExhibition.Range("F13").Value = 14476
$wsExhibition.Range("F21").Value = 44
$wsExhibition.Range("F25").Value = 5524
$wsExhibition.Range("F26").Value = 943
$wsExhibition.Range("F27").Value = 1034
$wsExhibition.Range("F28").Value = 352
$wsExhibition.Range("F30").Value = 114

# --- 全部类型 sheet ---
$wsAllTypes.Range("F3").Value  = 2223
$wsAllTypes.Range("F5").Value  = 13339
$wsAllTypes.Range("F9").Value  = 521
$wsAllTypes.Range("F11").Value = 1195
$wsAllTypes.Range("F12").Value = 1003
$wsAllTypes.Range("F14").Value = 14476
$wsAllTypes.Range("F22").Value = 44
$wsAllTypes.Range("F26").Value = 5524
$wsAllTypes.Range("F27").Value = 943
$wsAllTypes.Range("F28").Value = 1034
$wsAllTypes.Range("F29").Value = 352
$wsAllTypes.Range("F31").Value = 114
